$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 5, pushing existing rows 5-7 down to 6-8.
$ws.Rows.Item(5).Insert()

# New row 5 values
$ws.Range("A5").Value = 7
$ws.Range("B5").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C5").Value = "Ñuble"
$ws.Range("D5").Value = 44435
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
$ws.Range("E5").Value = 16
$ws.Range("F5").Value = 100112013
$ws.Range("G5").Value = "Alcachofa"
$ws.Range("H5").Value = "Madrigal"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 14000
$ws.Range("L5").Value = 15000
$ws.Range("M5").Value = 14500
$ws.Range("N5").Value = '$/caja 40 unidades'
$ws.Range("O5").Value = "Provincia del Elquí"
$ws.Range("P5").Value = 362
$ws.Range("Q5").Value = 40
$ws.Range("R5").Value = "Hortaliza"

# New row 9 (appended after current last row 8)
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C9").Value = "Ñuble"
$ws.Range("D9").Value = 44432
$ws.Range("D9").NumberFormat = $ws.Range("D8").NumberFormat
$ws.Range("E9").Value = 16
$ws.Range("F9").Value = 100112013
$ws.Range("G9").Value = "Alcachofa"
$ws.Range("H9").Value = "Madrigal"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 120
$ws.Range("K9").Value = 14000
$ws.Range("L9").Value = 15000
$ws.Range("M9").Value = 14500
$ws.Range("N9").Value = '$/caja 40 unidades'
$ws.Range("O9").Value = "Provincia del Elquí"
$ws.Range("P9").Value = 362
$ws.Range("Q9").Value = 40
$ws.Range("R9").Value = "Hortaliza"
